$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 8 (old rows 8+ shift down by one) ---
$ws.Rows.Item(8).Insert()

# --- Insert two new (blank) rows at 14:15 (pushes the summary/header block down) ---
$ws.Range("14:15").Insert()

# --- Row 7: task text changed to a brand-new string ---
$ws.Cells.Item(7,2).Value = "Макет страницы выбранной задачи"

# --- New row 8: restore the task that used to live in row 7 ---
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "Макет страницы выбора раздела"
$ws.Cells.Item(8,3).Value = "Петров"
$ws.Cells.Item(8,2).Style = $ws.Cells.Item(2,2).Style
$ws.Cells.Item(8,3).Style = $ws.Cells.Item(2,3).Style

# --- Renumber the "#" column for the rows pushed down by the insert ---
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(12,1).Value = 11

# --- Column D ("Степень выполнения") percentages for rows 2-11 ---
$ws.Cells.Item(2,4).Value = 1
$ws.Cells.Item(3,4).Value = 1
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(5,4).Value = 0.9
$ws.Cells.Item(6,4).Value = 1
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(8,4).Value = 1
$ws.Cells.Item(9,4).Value = 1
$ws.Cells.Item(10,4).Value = 0.5
$ws.Cells.Item(11,4).Value = 0

# --- Row 10 (was row 9): reviewer changed from "Петров" to "Руданов" ---
$ws.Cells.Item(10,3).Value = "Руданов"

# --- Row 12 D cell: drop the old "," placeholder value ---
$ws.Cells.Item(12,4).ClearContents()

# --- Materialize the empty C14 placeholder cell (copied from the row above) ---
$ws.Cells.Item(13,3).Copy($ws.Cells.Item(14,3))

# --- Clean up stray cells left behind by the row-insert copy ---
$ws.Cells.Item(16,3).ClearContents()
$ws.Cells.Item(18,4).Clear()
$ws.Cells.Item(19,4).Clear()

# --- Selection, as recorded after the edit ---
$ws.Range("D8").Select()
